# Update the quarterly income-statement data: drop oldest quarter (1399/06),
# shift remaining quarters left, append new quarter (1401/12) with updated
# figures from the revised read_price algorithm.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period headers (D8:M8) ---
$ws.Cells.Item(8, 4).Value = "فصل سوم منتهی به 1399/09"
$ws.Cells.Item(8, 5).Value = "فصل چهارم منتهی به 1399/12"
$ws.Cells.Item(8, 6).Value = "فصل اول منتهی به 1400/03"
$ws.Cells.Item(8, 7).Value = "فصل دوم منتهی به 1400/06"
$ws.Cells.Item(8, 8).Value = "فصل سوم منتهی به 1400/09"
$ws.Cells.Item(8, 9).Value = "فصل چهارم منتهی به 1400/12"
$ws.Cells.Item(8, 10).Value = "فصل اول منتهی به 1401/03"
$ws.Cells.Item(8, 11).Value = "فصل دوم منتهی به 1401/06"
$ws.Cells.Item(8, 12).Value = "فصل سوم منتهی به 1401/09"
$ws.Cells.Item(8, 13).Value = "فصل چهارم منتهی به 1401/12"

# --- Row 9: publish dates (D9:M9) ---
$ws.Cells.Item(9, 4).Value = "1400-10-29 (3)"
$ws.Cells.Item(9, 5).Value = "1401-02-10 (12)"
$ws.Cells.Item(9, 6).Value = "1401-05-12 (4)"
$ws.Cells.Item(9, 7).Value = "1401-09-09 (4)"
$ws.Cells.Item(9, 8).Value = "1401-10-29 (3)"
$ws.Cells.Item(9, 9).Value = "1402-02-09 (10)"
$ws.Cells.Item(9, 10).Value = "1401-05-12 (2)"
$ws.Cells.Item(9, 11).Value = "1401-09-09 (2)"
$ws.Cells.Item(9, 12).Value = "1401-10-29"
$ws.Cells.Item(9, 13).Value = "1402-02-09 (2)"

# --- Row 11: Sales (فروش) (D11:M11) ---
$ws.Cells.Item(11, 4).Value = 1958
$ws.Cells.Item(11, 5).Value = 2504
$ws.Cells.Item(11, 6).Value = 3630
$ws.Cells.Item(11, 7).Value = 2646
$ws.Cells.Item(11, 8).Value = 1991
$ws.Cells.Item(11, 9).Value = 3463
$ws.Cells.Item(11, 10).Value = 3819
$ws.Cells.Item(11, 11).Value = 2693
$ws.Cells.Item(11, 12).Value = 3303
$ws.Cells.Item(11, 13).Value = 2427

# --- Row 12: COGS (بهای تمام شده کالای فروش رفته) (D12:M12) ---
$ws.Cells.Item(12, 4).Value = -998
$ws.Cells.Item(12, 5).Value = -1509
$ws.Cells.Item(12, 6).Value = -1969
$ws.Cells.Item(12, 7).Value = -1329
$ws.Cells.Item(12, 8).Value = -1339
$ws.Cells.Item(12, 9).Value = -2119
$ws.Cells.Item(12, 10).Value = -2409
$ws.Cells.Item(12, 11).Value = -1977
$ws.Cells.Item(12, 12).Value = -1975
$ws.Cells.Item(12, 13).Value = -1879

# --- Row 13: Gross profit (سود ناخالص) (D13:M13) ---
$ws.Cells.Item(13, 4).Value = 960
$ws.Cells.Item(13, 5).Value = 995
$ws.Cells.Item(13, 6).Value = 1661
$ws.Cells.Item(13, 7).Value = 1316
$ws.Cells.Item(13, 8).Value = 652
$ws.Cells.Item(13, 9).Value = 1345
$ws.Cells.Item(13, 10).Value = 1410
$ws.Cells.Item(13, 11).Value = 716
$ws.Cells.Item(13, 12).Value = 1328
$ws.Cells.Item(13, 13).Value = 548

# --- Row 14: SG&A expenses (D14:M14) ---
$ws.Cells.Item(14, 4).Value = -89
$ws.Cells.Item(14, 5).Value = -106
$ws.Cells.Item(14, 6).Value = -165
$ws.Cells.Item(14, 7).Value = -150
$ws.Cells.Item(14, 8).Value = -127
$ws.Cells.Item(14, 9).Value = -155
$ws.Cells.Item(14, 10).Value = -194
$ws.Cells.Item(14, 11).Value = -160
$ws.Cells.Item(14, 12).Value = -266
$ws.Cells.Item(14, 13).Value = -142

# --- Row 15: Impairment expense (D15:M15) ---
$ws.Cells.Item(15, 4).Value = "-"
$ws.Cells.Item(15, 5).Value = 10
$ws.Cells.Item(15, 6).Value = "-"
$ws.Cells.Item(15, 7).Value = "-"
$ws.Cells.Item(15, 8).Value = "-"
$ws.Cells.Item(15, 9).Value = "-"
$ws.Cells.Item(15, 10).Value = "-"
$ws.Cells.Item(15, 11).Value = "-"
$ws.Cells.Item(15, 12).Value = "-"
$ws.Cells.Item(15, 13).Value = "-"

# --- Row 16: Other operating income (expenses), net (D16:M16) ---
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 5).Value = -1
$ws.Cells.Item(16, 6).Value = 22
$ws.Cells.Item(16, 7).Value = -17
$ws.Cells.Item(16, 8).Value = 170
$ws.Cells.Item(16, 9).Value = -178
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = -17
$ws.Cells.Item(16, 12).Value = 2
$ws.Cells.Item(16, 13).Value = 2

# --- Row 17: Operating profit (D17:M17) ---
$ws.Cells.Item(17, 4).Value = 872
$ws.Cells.Item(17, 5).Value = 898
$ws.Cells.Item(17, 6).Value = 1517
$ws.Cells.Item(17, 7).Value = 1149
$ws.Cells.Item(17, 8).Value = 695
$ws.Cells.Item(17, 9).Value = 1012
$ws.Cells.Item(17, 10).Value = 1217
$ws.Cells.Item(17, 11).Value = 539
$ws.Cells.Item(17, 12).Value = 1065
$ws.Cells.Item(17, 13).Value = 407

# --- Row 18: Finance costs (D18:M18) ---
$ws.Cells.Item(18, 4).Value = -2
$ws.Cells.Item(18, 5).Value = -3
$ws.Cells.Item(18, 6).Value = "-"
$ws.Cells.Item(18, 7).Value = "-"
$ws.Cells.Item(18, 8).Value = "-"
$ws.Cells.Item(18, 9).Value = "-"
$ws.Cells.Item(18, 10).Value = "-"
$ws.Cells.Item(18, 11).Value = -3
$ws.Cells.Item(18, 12).Value = -7
$ws.Cells.Item(18, 13).Value = -6

# --- Row 19: Other non-operating income (expenses), net (D19:M19) ---
$ws.Cells.Item(19, 4).Value = 101
$ws.Cells.Item(19, 5).Value = 82
$ws.Cells.Item(19, 6).Value = 33
$ws.Cells.Item(19, 7).Value = 40
$ws.Cells.Item(19, 8).Value = 67
$ws.Cells.Item(19, 9).Value = -167
$ws.Cells.Item(19, 10).Value = 6
$ws.Cells.Item(19, 11).Value = 31
$ws.Cells.Item(19, 12).Value = -6
$ws.Cells.Item(19, 13).Value = 14

# --- Row 20: Profit before tax from continuing operations (D20:M20) ---
$ws.Cells.Item(20, 4).Value = 971
$ws.Cells.Item(20, 5).Value = 976
$ws.Cells.Item(20, 6).Value = 1550
$ws.Cells.Item(20, 7).Value = 1189
$ws.Cells.Item(20, 8).Value = 761
$ws.Cells.Item(20, 9).Value = 845
$ws.Cells.Item(20, 10).Value = 1222
$ws.Cells.Item(20, 11).Value = 567
$ws.Cells.Item(20, 12).Value = 1051
$ws.Cells.Item(20, 13).Value = 415

# --- Row 21: Tax (D21:M21) ---
$ws.Cells.Item(21, 4).Value = -174
$ws.Cells.Item(21, 5).Value = 45
$ws.Cells.Item(21, 6).Value = -303
$ws.Cells.Item(21, 7).Value = -230
$ws.Cells.Item(21, 8).Value = -139
$ws.Cells.Item(21, 9).Value = -43
$ws.Cells.Item(21, 10).Value = -183
$ws.Cells.Item(21, 11).Value = 492
$ws.Cells.Item(21, 12).Value = -307
$ws.Cells.Item(21, 13).Value = -100

# --- Row 22: Profit from continuing operations (D22:M22) ---
$ws.Cells.Item(22, 4).Value = 797
$ws.Cells.Item(22, 5).Value = 1021
$ws.Cells.Item(22, 6).Value = 1247
$ws.Cells.Item(22, 7).Value = 959
$ws.Cells.Item(22, 8).Value = 622
$ws.Cells.Item(22, 9).Value = 802
$ws.Cells.Item(22, 10).Value = 1039
$ws.Cells.Item(22, 11).Value = 1059
$ws.Cells.Item(22, 12).Value = 744
$ws.Cells.Item(22, 13).Value = 316

# --- Row 23: Discontinued operations (D23:M23) ---
$ws.Cells.Item(23, 4).Value = "-"
$ws.Cells.Item(23, 5).Value = "-"
$ws.Cells.Item(23, 6).Value = "-"
$ws.Cells.Item(23, 7).Value = "-"
$ws.Cells.Item(23, 8).Value = "-"
$ws.Cells.Item(23, 9).Value = "-"
$ws.Cells.Item(23, 10).Value = "-"
$ws.Cells.Item(23, 11).Value = "-"
$ws.Cells.Item(23, 12).Value = "-"
$ws.Cells.Item(23, 13).Value = "-"

# --- Row 24: Net profit (D24:M24) ---
$ws.Cells.Item(24, 4).Value = 797
$ws.Cells.Item(24, 5).Value = 1021
$ws.Cells.Item(24, 6).Value = 1247
$ws.Cells.Item(24, 7).Value = 959
$ws.Cells.Item(24, 8).Value = 622
$ws.Cells.Item(24, 9).Value = 802
$ws.Cells.Item(24, 10).Value = 1039
$ws.Cells.Item(24, 11).Value = 1059
$ws.Cells.Item(24, 12).Value = 744
$ws.Cells.Item(24, 13).Value = 316

# --- Row 25: EPS after tax (D25:M25) ---
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 0

# --- Row 26: Capital (D26:M26) ---
$ws.Cells.Item(26, 4).Value = 1656
$ws.Cells.Item(26, 5).Value = 3292
$ws.Cells.Item(26, 6).Value = 3446
$ws.Cells.Item(26, 7).Value = 3080
$ws.Cells.Item(26, 8).Value = 2824
$ws.Cells.Item(26, 9).Value = 7261
$ws.Cells.Item(26, 10).Value = 2736
$ws.Cells.Item(26, 11).Value = 6448
$ws.Cells.Item(26, 12).Value = 5757
$ws.Cells.Item(26, 13).Value = 4400

# --- Row 27: EPS on latest capital (D27:M27) ---
$ws.Cells.Item(27, 4).Value = 0
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 13).Value = 0

# --- Column widths: the width pattern (29,29,29,31 repeating) shifts left by
# one column along with the data so the "wide" (31) column stays aligned on
# the same relative quarter position. ---
$ws.Columns.Item(4).ColumnWidth = 28.17
$ws.Columns.Item(5).ColumnWidth = 30.17
$ws.Columns.Item(6).ColumnWidth = 28.17
$ws.Columns.Item(7).ColumnWidth = 28.17
$ws.Columns.Item(8).ColumnWidth = 28.17
$ws.Columns.Item(9).ColumnWidth = 30.17
$ws.Columns.Item(10).ColumnWidth = 28.17
$ws.Columns.Item(11).ColumnWidth = 28.17
$ws.Columns.Item(12).ColumnWidth = 28.17
$ws.Columns.Item(13).ColumnWidth = 30.17
